$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 2).Value = 6139017
$ws.Cells.Item(4, 5).Value = 'JK Tammeka Tartu'
$ws.Cells.Item(4, 6).Value = 'Harju JK Laagri'
$ws.Cells.Item(4, 7).Value = 2
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 'H'
$ws.Cells.Item(4, 10).Value = 1.666
$ws.Cells.Item(4, 11).Value = 3.6
$ws.Cells.Item(4, 12).Value = 4.2
$ws.Cells.Item(4, 13).Value = 1.727
$ws.Cells.Item(4, 14).Value = 3.5
$ws.Cells.Item(4, 15).Value = 4
$ws.Cells.Item(4, 16).Value = -0.75
$ws.Cells.Item(4, 17).Value = 2
$ws.Cells.Item(4, 18).Value = 1.8
$ws.Cells.Item(4, 19).Value = 2.5
$ws.Cells.Item(4, 20).Value = 1.9
$ws.Cells.Item(4, 21).Value = 1.9
$ws.Cells.Item(4, 22).Value = 0.7270000000000001
$ws.Cells.Item(4, 24).Value = -1
$ws.Cells.Item(4, 25).Value = 1
$ws.Cells.Item(4, 26).Value = -1
$ws.Cells.Item(4, 28).Value = 0.8999999999999999

# Row 5
$ws.Cells.Item(5, 2).Value = 6139018
$ws.Cells.Item(5, 5).Value = 'JK Tallinna Kalev'
$ws.Cells.Item(5, 6).Value = 'JK Trans Narva'
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 1
$ws.Cells.Item(5, 9).Value = 'A'
$ws.Cells.Item(5, 10).Value = 2.4
$ws.Cells.Item(5, 11).Value = 3.4
$ws.Cells.Item(5, 12).Value = 2.5
$ws.Cells.Item(5, 13).Value = 2.875
$ws.Cells.Item(5, 14).Value = 3.1
$ws.Cells.Item(5, 15).Value = 2.3
$ws.Cells.Item(5, 16).Value = 0.25
$ws.Cells.Item(5, 17).Value = 1.75
$ws.Cells.Item(5, 18).Value = 2.05
$ws.Cells.Item(5, 19).Value = 2.25
$ws.Cells.Item(5, 20).Value = 1.925
$ws.Cells.Item(5, 21).Value = 1.875
$ws.Cells.Item(5, 22).Value = -1
$ws.Cells.Item(5, 24).Value = 1.3
$ws.Cells.Item(5, 25).Value = -1
$ws.Cells.Item(5, 26).Value = 1.05
$ws.Cells.Item(5, 28).Value = 0.875

# Row 10
$ws.Cells.Item(10, 6).Value = 'JK Tammeka Tartu'

# Row 11
$ws.Cells.Item(11, 6).Value = 'JK Tallinna Kalev'

# Row 13
$ws.Cells.Item(13, 5).Value = 'JK Tallinna Kalev'

# Row 15
$ws.Cells.Item(15, 5).Value = 'JK Tammeka Tartu'

# Row 20
$ws.Cells.Item(20, 6).Value = 'JK Tammeka Tartu'

# Row 21
$ws.Cells.Item(21, 6).Value = 'JK Tallinna Kalev'

# Row 24
$ws.Cells.Item(24, 5).Value = 'JK Tammeka Tartu'

# Row 25
$ws.Cells.Item(25, 6).Value = 'JK Tallinna Kalev'

# Row 27
$ws.Cells.Item(27, 5).Value = 'JK Tammeka Tartu'

# Row 30
$ws.Cells.Item(30, 5).Value = 'JK Tallinna Kalev'

# Row 33
$ws.Cells.Item(33, 6).Value = 'JK Tallinna Kalev'

# Row 36
$ws.Cells.Item(36, 5).Value = 'JK Tammeka Tartu'
$ws.Cells.Item(36, 6).Value = 'JK Tallinna Kalev'

# Row 39
$ws.Cells.Item(39, 6).Value = 'JK Tammeka Tartu'

# Row 41
$ws.Cells.Item(41, 5).Value = 'JK Tallinna Kalev'

# Row 44
$ws.Cells.Item(44, 6).Value = 'JK Tallinna Kalev'

# Row 47
$ws.Cells.Item(47, 5).Value = 'JK Tammeka Tartu'

# Row 49
$ws.Cells.Item(49, 5).Value = 'JK Tallinna Kalev'

# Row 50
$ws.Cells.Item(50, 6).Value = 'JK Tammeka Tartu'

# Row 53
$ws.Cells.Item(53, 6).Value = 'JK Tammeka Tartu'

# Row 55
$ws.Cells.Item(55, 6).Value = 'JK Tallinna Kalev'

# Row 58
$ws.Cells.Item(58, 5).Value = 'JK Tallinna Kalev'

# Row 59
$ws.Cells.Item(59, 5).Value = 'JK Tammeka Tartu'

# Row 63
$ws.Cells.Item(63, 6).Value = 'JK Tammeka Tartu'

# Row 64
$ws.Cells.Item(64, 2).Value = 6139067
$ws.Cells.Item(64, 5).Value = 'Paide Linnameeskond'
$ws.Cells.Item(64, 6).Value = 'Parnu JK Vaprus'
$ws.Cells.Item(64, 7).Value = 3
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 'H'
$ws.Cells.Item(64, 10).Value = 1.8
$ws.Cells.Item(64, 11).Value = 3.4
$ws.Cells.Item(64, 13).Value = 1.5
$ws.Cells.Item(64, 14).Value = 3.8
$ws.Cells.Item(64, 17).Value = 1.75
$ws.Cells.Item(64, 22).Value = 0.5
$ws.Cells.Item(64, 24).Value = -1
$ws.Cells.Item(64, 25).Value = 0.75
$ws.Cells.Item(64, 26).Value = -1

# Row 65
$ws.Cells.Item(65, 2).Value = 6139064
$ws.Cells.Item(65, 5).Value = 'JK Trans Narva'
$ws.Cells.Item(65, 6).Value = 'Harju JK Laagri'
$ws.Cells.Item(65, 7).Value = 1
$ws.Cells.Item(65, 8).Value = 3
$ws.Cells.Item(65, 9).Value = 'A'
$ws.Cells.Item(65, 10).Value = 1.75
$ws.Cells.Item(65, 11).Value = 3.6
$ws.Cells.Item(65, 13).Value = 1.45
$ws.Cells.Item(65, 14).Value = 4
$ws.Cells.Item(65, 17).Value = 1.85
$ws.Cells.Item(65, 22).Value = -1
$ws.Cells.Item(65, 24).Value = 5
$ws.Cells.Item(65, 25).Value = -1
$ws.Cells.Item(65, 26).Value = 0.95

# Row 66
$ws.Cells.Item(66, 6).Value = 'JK Tallinna Kalev'

# Row 70
$ws.Cells.Item(70, 5).Value = 'JK Tammeka Tartu'

# Row 71
$ws.Cells.Item(71, 2).Value = 6139071
$ws.Cells.Item(71, 5).Value = 'Parnu JK Vaprus'
$ws.Cells.Item(71, 6).Value = 'JK Trans Narva'
$ws.Cells.Item(71, 7).Value = 3
$ws.Cells.Item(71, 9).Value = 'H'
$ws.Cells.Item(71, 10).Value = 2.4
$ws.Cells.Item(71, 11).Value = 3.2
$ws.Cells.Item(71, 12).Value = 2.6
$ws.Cells.Item(71, 13).Value = 3
$ws.Cells.Item(71, 14).Value = 3.25
$ws.Cells.Item(71, 15).Value = 2.2
$ws.Cells.Item(71, 16).Value = 0.25
$ws.Cells.Item(71, 17).Value = 1.825
$ws.Cells.Item(71, 18).Value = 1.975
$ws.Cells.Item(71, 19).Value = 2.5
$ws.Cells.Item(71, 20).Value = 1.875
$ws.Cells.Item(71, 21).Value = 1.925
$ws.Cells.Item(71, 22).Value = 2
$ws.Cells.Item(71, 24).Value = -1
$ws.Cells.Item(71, 25).Value = 0.825
$ws.Cells.Item(71, 27).Value = 0.875
$ws.Cells.Item(71, 28).Value = -1

# Row 72
$ws.Cells.Item(72, 2).Value = 6139072
$ws.Cells.Item(72, 5).Value = 'JK Tammeka Tartu'
$ws.Cells.Item(72, 6).Value = 'FC Flora Tallinn'
$ws.Cells.Item(72, 7).Value = 1
$ws.Cells.Item(72, 9).Value = 'A'
$ws.Cells.Item(72, 10).Value = 9
$ws.Cells.Item(72, 11).Value = 7
$ws.Cells.Item(72, 12).Value = 1.166
$ws.Cells.Item(72, 13).Value = 7
$ws.Cells.Item(72, 14).Value = 6
$ws.Cells.Item(72, 15).Value = 1.25
$ws.Cells.Item(72, 16).Value = 1.75
$ws.Cells.Item(72, 17).Value = 1.9
$ws.Cells.Item(72, 18).Value = 1.9
$ws.Cells.Item(72, 19).Value = 3
$ws.Cells.Item(72, 20).Value = 1.95
$ws.Cells.Item(72, 21).Value = 1.85
$ws.Cells.Item(72, 22).Value = -1
$ws.Cells.Item(72, 24).Value = 0.25
$ws.Cells.Item(72, 25).Value = 0.8999999999999999
$ws.Cells.Item(72, 27).Value = 0
$ws.Cells.Item(72, 28).Value = 0

# Row 74
$ws.Cells.Item(74, 5).Value = 'JK Tallinna Kalev'

# Row 75
$ws.Cells.Item(75, 6).Value = 'JK Tammeka Tartu'

# Row 79
$ws.Cells.Item(79, 5).Value = 'JK Tammeka Tartu'

# Row 80
$ws.Cells.Item(80, 6).Value = 'JK Tallinna Kalev'

# Row 83
$ws.Cells.Item(83, 5).Value = 'JK Tammeka Tartu'

# Row 85
$ws.Cells.Item(85, 6).Value = 'JK Tallinna Kalev'

# Row 88
$ws.Cells.Item(88, 2).Value = 6376945
$ws.Cells.Item(88, 5).Value = 'Parnu JK Vaprus'
$ws.Cells.Item(88, 6).Value = 'Harju JK Laagri'
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 'D'
$ws.Cells.Item(88, 10).Value = 1.615
$ws.Cells.Item(88, 11).Value = 4
$ws.Cells.Item(88, 12).Value = 4.5
$ws.Cells.Item(88, 13).Value = 1.85
$ws.Cells.Item(88, 14).Value = 3.8
$ws.Cells.Item(88, 15).Value = 3.5
$ws.Cells.Item(88, 16).Value = -0.5
$ws.Cells.Item(88, 19).Value = 2.5
$ws.Cells.Item(88, 20).Value = 1.75
$ws.Cells.Item(88, 21).Value = 1.95
$ws.Cells.Item(88, 23).Value = 2.8
$ws.Cells.Item(88, 24).Value = -1
$ws.Cells.Item(88, 27).Value = -1
$ws.Cells.Item(88, 28).Value = 0.95

# Row 89
$ws.Cells.Item(89, 2).Value = 6376947
$ws.Cells.Item(89, 5).Value = 'JK Tammeka Tartu'
$ws.Cells.Item(89, 6).Value = 'JK Tallinna Kalev'
$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(89, 8).Value = 7
$ws.Cells.Item(89, 9).Value = 'A'
$ws.Cells.Item(89, 10).Value = 3.6
$ws.Cells.Item(89, 11).Value = 3.4
$ws.Cells.Item(89, 12).Value = 1.909
$ws.Cells.Item(89, 13).Value = 2.4
$ws.Cells.Item(89, 14).Value = 3.6
$ws.Cells.Item(89, 15).Value = 2.45
$ws.Cells.Item(89, 16).Value = 0
$ws.Cells.Item(89, 19).Value = 2.75
$ws.Cells.Item(89, 20).Value = 1.975
$ws.Cells.Item(89, 21).Value = 1.825
$ws.Cells.Item(89, 23).Value = -1
$ws.Cells.Item(89, 24).Value = 1.45
$ws.Cells.Item(89, 27).Value = 0.9750000000000001
$ws.Cells.Item(89, 28).Value = -1

# Row 93
$ws.Cells.Item(93, 6).Value = 'JK Tammeka Tartu'

# Row 96
$ws.Cells.Item(96, 5).Value = 'JK Tallinna Kalev'

# Row 100
$ws.Cells.Item(100, 5).Value = 'JK Tammeka Tartu'

# Row 102
$ws.Cells.Item(102, 6).Value = 'JK Tallinna Kalev'

# Row 103
$ws.Cells.Item(103, 5).Value = 'JK Tallinna Kalev'

# Row 104
$ws.Cells.Item(104, 2).Value = 6535416
$ws.Cells.Item(104, 5).Value = 'Paide Linnameeskond'
$ws.Cells.Item(104, 6).Value = 'FC Levadia Tallinn'
$ws.Cells.Item(104, 7).Value = 2
$ws.Cells.Item(104, 8).Value = 2
$ws.Cells.Item(104, 9).Value = 'D'
$ws.Cells.Item(104, 10).Value = 3
$ws.Cells.Item(104, 11).Value = 3.8
$ws.Cells.Item(104, 12).Value = 2
$ws.Cells.Item(104, 13).Value = 3
$ws.Cells.Item(104, 14).Value = 4
$ws.Cells.Item(104, 15).Value = 1.909
$ws.Cells.Item(104, 16).Value = 0.5
$ws.Cells.Item(104, 17).Value = 1.85
$ws.Cells.Item(104, 18).Value = 1.95
$ws.Cells.Item(104, 22).Value = -1
$ws.Cells.Item(104, 23).Value = 3
$ws.Cells.Item(104, 25).Value = 0.8500000000000001
$ws.Cells.Item(104, 27).Value = 0.95
$ws.Cells.Item(104, 28).Value = -1

# Row 105
$ws.Cells.Item(105, 2).Value = 6533597
$ws.Cells.Item(105, 5).Value = 'FC Kuressaare'
$ws.Cells.Item(105, 6).Value = 'Parnu JK Vaprus'
$ws.Cells.Item(105, 7).Value = 1
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 9).Value = 'H'
$ws.Cells.Item(105, 10).Value = 2.5
$ws.Cells.Item(105, 11).Value = 3.4
$ws.Cells.Item(105, 12).Value = 2.5
$ws.Cells.Item(105, 13).Value = 2.15
$ws.Cells.Item(105, 14).Value = 3.6
$ws.Cells.Item(105, 15).Value = 2.875
$ws.Cells.Item(105, 16).Value = -0.25
$ws.Cells.Item(105, 17).Value = 1.95
$ws.Cells.Item(105, 18).Value = 1.85
$ws.Cells.Item(105, 22).Value = 1.15
$ws.Cells.Item(105, 23).Value = -1
$ws.Cells.Item(105, 25).Value = 0.95
$ws.Cells.Item(105, 27).Value = -1
$ws.Cells.Item(105, 28).Value = 0.8500000000000001

# Row 106
$ws.Cells.Item(106, 2).Value = 6537957
$ws.Cells.Item(106, 5).Value = 'FC Flora Tallinn'
$ws.Cells.Item(106, 6).Value = 'JK Nomme Kalju'
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 9).Value = 'D'
$ws.Cells.Item(106, 10).Value = 1.4
$ws.Cells.Item(106, 12).Value = 7.5
$ws.Cells.Item(106, 13).Value = 1.5
$ws.Cells.Item(106, 14).Value = 4.2
$ws.Cells.Item(106, 15).Value = 5
$ws.Cells.Item(106, 16).Value = -1
$ws.Cells.Item(106, 17).Value = 1.85
$ws.Cells.Item(106, 18).Value = 1.95
$ws.Cells.Item(106, 20).Value = 1.85
$ws.Cells.Item(106, 21).Value = 1.95
$ws.Cells.Item(106, 22).Value = -1
$ws.Cells.Item(106, 23).Value = 3.2
$ws.Cells.Item(106, 25).Value = -1
$ws.Cells.Item(106, 26).Value = 0.95
$ws.Cells.Item(106, 27).Value = -1
$ws.Cells.Item(106, 28).Value = 0.95

# Row 107
$ws.Cells.Item(107, 2).Value = 6537869
$ws.Cells.Item(107, 5).Value = 'JK Tallinna Kalev'
$ws.Cells.Item(107, 6).Value = 'JK Trans Narva'
$ws.Cells.Item(107, 7).Value = 5
$ws.Cells.Item(107, 9).Value = 'H'
$ws.Cells.Item(107, 10).Value = 1.6
$ws.Cells.Item(107, 12).Value = 4.5
$ws.Cells.Item(107, 13).Value = 1.65
$ws.Cells.Item(107, 14).Value = 4
$ws.Cells.Item(107, 15).Value = 4.333
$ws.Cells.Item(107, 16).Value = -0.75
$ws.Cells.Item(107, 17).Value = 1.8
$ws.Cells.Item(107, 18).Value = 2
$ws.Cells.Item(107, 20).Value = 1.9
$ws.Cells.Item(107, 21).Value = 1.9
$ws.Cells.Item(107, 22).Value = 0.6499999999999999
$ws.Cells.Item(107, 23).Value = -1
$ws.Cells.Item(107, 25).Value = 0.8
$ws.Cells.Item(107, 26).Value = -1
$ws.Cells.Item(107, 27).Value = 0.8999999999999999
$ws.Cells.Item(107, 28).Value = -1

# Row 108
$ws.Cells.Item(108, 5).Value = 'JK Tallinna Kalev'

# Row 112
$ws.Cells.Item(112, 5).Value = 'JK Tammeka Tartu'
$ws.Cells.Item(112, 6).Value = 'JK Tallinna Kalev'

# Row 114
$ws.Cells.Item(114, 5).Value = 'JK Tallinna Kalev'

# Row 115
$ws.Cells.Item(115, 2).Value = 7919322
$ws.Cells.Item(115, 5).Value = 'FC Kuressaare'
$ws.Cells.Item(115, 6).Value = 'FC Levadia Tallinn'
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 6
$ws.Cells.Item(115, 9).Value = 'A'
$ws.Cells.Item(115, 10).Value = 11
$ws.Cells.Item(115, 11).Value = 6
$ws.Cells.Item(115, 12).Value = 1.166
$ws.Cells.Item(115, 13).Value = 15
$ws.Cells.Item(115, 14).Value = 8.5
$ws.Cells.Item(115, 15).Value = 1.125
$ws.Cells.Item(115, 16).Value = 2.5
$ws.Cells.Item(115, 17).Value = 1.825
$ws.Cells.Item(115, 18).Value = 1.975
$ws.Cells.Item(115, 19).Value = 3.25
$ws.Cells.Item(115, 20).Value = 1.9
$ws.Cells.Item(115, 21).Value = 1.9
$ws.Cells.Item(115, 22).Value = -1
$ws.Cells.Item(115, 24).Value = 0.125
$ws.Cells.Item(115, 25).Value = -1
$ws.Cells.Item(115, 26).Value = 0.9750000000000001
$ws.Cells.Item(115, 27).Value = 0.8999999999999999
$ws.Cells.Item(115, 28).Value = -1

# Row 116
$ws.Cells.Item(116, 2).Value = 7919323
$ws.Cells.Item(116, 5).Value = 'JK Nomme Kalju'
$ws.Cells.Item(116, 6).Value = 'JK Trans Narva'
$ws.Cells.Item(116, 7).Value = 3
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 9).Value = 'H'
$ws.Cells.Item(116, 10).Value = 1.285
$ws.Cells.Item(116, 11).Value = 5.5
$ws.Cells.Item(116, 12).Value = 6.5
$ws.Cells.Item(116, 13).Value = 1.571
$ws.Cells.Item(116, 14).Value = 4.75
$ws.Cells.Item(116, 15).Value = 4.2
$ws.Cells.Item(116, 16).Value = -1
$ws.Cells.Item(116, 17).Value = 1.925
$ws.Cells.Item(116, 18).Value = 1.875
$ws.Cells.Item(116, 19).Value = 2.75
$ws.Cells.Item(116, 20).Value = 1.875
$ws.Cells.Item(116, 21).Value = 1.925
$ws.Cells.Item(116, 22).Value = 0.571
$ws.Cells.Item(116, 24).Value = -1
$ws.Cells.Item(116, 25).Value = 0.925
$ws.Cells.Item(116, 26).Value = -1
$ws.Cells.Item(116, 27).Value = 0.4375
$ws.Cells.Item(116, 28).Value = -0.5

# Row 118
$ws.Cells.Item(118, 6).Value = 'JK Tammeka Tartu'

# Row 119
$ws.Cells.Item(119, 5).Value = 'JK Tammeka Tartu'

# Row 120
$ws.Cells.Item(120, 2).Value = 7721007
$ws.Cells.Item(120, 5).Value = 'JK Trans Narva'
$ws.Cells.Item(120, 6).Value = 'JK Tammeka Tartu'
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 5
$ws.Cells.Item(120, 9).Value = 'A'
$ws.Cells.Item(120, 10).Value = 2.25
$ws.Cells.Item(120, 12).Value = 2.75
$ws.Cells.Item(120, 13).Value = 2.1
$ws.Cells.Item(120, 14).Value = 3.25
$ws.Cells.Item(120, 15).Value = 3
$ws.Cells.Item(120, 16).Value = -0.25
$ws.Cells.Item(120, 17).Value = 1.875
$ws.Cells.Item(120, 18).Value = 1.925
$ws.Cells.Item(120, 20).Value = 1.825
$ws.Cells.Item(120, 21).Value = 1.975
$ws.Cells.Item(120, 22).Value = -1
$ws.Cells.Item(120, 24).Value = 2
$ws.Cells.Item(120, 25).Value = -1
$ws.Cells.Item(120, 26).Value = 0.925
$ws.Cells.Item(120, 27).Value = 0.825

# Row 121
$ws.Cells.Item(121, 2).Value = 7721087
$ws.Cells.Item(121, 5).Value = 'Paide Linnameeskond'
$ws.Cells.Item(121, 6).Value = 'FC Flora Tallinn'
$ws.Cells.Item(121, 7).Value = 2
$ws.Cells.Item(121, 8).Value = 1
$ws.Cells.Item(121, 9).Value = 'H'
$ws.Cells.Item(121, 10).Value = 2.2
$ws.Cells.Item(121, 12).Value = 2.8
$ws.Cells.Item(121, 13).Value = 1.85
$ws.Cells.Item(121, 14).Value = 3.6
$ws.Cells.Item(121, 15).Value = 3.4
$ws.Cells.Item(121, 16).Value = -0.5
$ws.Cells.Item(121, 17).Value = 1.9
$ws.Cells.Item(121, 18).Value = 1.9
$ws.Cells.Item(121, 20).Value = 1.95
$ws.Cells.Item(121, 21).Value = 1.85
$ws.Cells.Item(121, 22).Value = 0.8500000000000001
$ws.Cells.Item(121, 24).Value = -1
$ws.Cells.Item(121, 25).Value = 0.8999999999999999
$ws.Cells.Item(121, 26).Value = -1
$ws.Cells.Item(121, 27).Value = 0.95

# Row 124
$ws.Cells.Item(124, 6).Value = 'JK Tammeka Tartu'

# Row 127
$ws.Cells.Item(127, 5).Value = 'JK Tallinna Kalev'

# Row 129
$ws.Cells.Item(129, 6).Value = 'JK Tallinna Kalev'

# Row 133
$ws.Cells.Item(133, 5).Value = 'JK Tammeka Tartu'

# Row 134
$ws.Cells.Item(134, 5).Value = 'JK Tallinna Kalev'

# Row 135
$ws.Cells.Item(135, 6).Value = 'JK Tammeka Tartu'

# Row 139
$ws.Cells.Item(139, 5).Value = 'JK Tammeka Tartu'

# Row 141
$ws.Cells.Item(141, 6).Value = 'JK Tammeka Tartu'

# Row 142
$ws.Cells.Item(142, 6).Value = 'JK Tallinna Kalev'

# Row 146
$ws.Cells.Item(146, 2).Value = '7721021'
$ws.Cells.Item(146, 4).Value = 45410.35416666666
$ws.Cells.Item(146, 5).Value = 'JK Tammeka Tartu'
$ws.Cells.Item(146, 6).Value = 'FC Levadia Tallinn'
$ws.Cells.Item(146, 10).Value = 7
$ws.Cells.Item(146, 11).Value = 6
$ws.Cells.Item(146, 12).Value = 1.25
$ws.Cells.Item(146, 13).Value = 7
$ws.Cells.Item(146, 14).Value = 5.75
$ws.Cells.Item(146, 15).Value = 1.25
$ws.Cells.Item(146, 16).Value = 1.75
$ws.Cells.Item(146, 17).Value = 1.85
$ws.Cells.Item(146, 18).Value = 1.95
$ws.Cells.Item(146, 19).Value = 3
